$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1958.8077
$ws.Range("J17").Value = 1958.8077
$ws.Range("L17").Value = 5876.4231
$ws.Range("N17").Value = -6212.4231

$ws.Range("H86").Value = 4786.143
$ws.Range("I86").Value = 1833.3334
$ws.Range("K86").Value = 1833.3334
$ws.Range("M86").Value = -710.3334

$ws.Range("H89").Value = 4786.143
$ws.Range("I89").Value = 1833.3334
$ws.Range("K89").Value = 9166.666999999999
$ws.Range("M89").Value = -3550.666999999999

$ws.Range("H106").Value = 2174.75
$ws.Range("J106").Value = 3200
$ws.Range("L106").Value = 3200
$ws.Range("N106").Value = -4462

$ws.Range("H116").Value = 34246.125
$ws.Range("I116").Value = 45259.812
$ws.Range("K116").Value = 45259.812
$ws.Range("M116").Value = -41817.812

$ws.Range("H137").Value = 2218
$ws.Range("I137").Value = 1170
$ws.Range("K137").Value = 3510
$ws.Range("M137").Value = -960

$ws.Range("H138").Value = 2596
$ws.Range("I138").Value = 1655.3572
$ws.Range("J138").Value = 3793.182
$ws.Range("K138").Value = 4966.071599999999
$ws.Range("L138").Value = 11379.546
$ws.Range("M138").Value = 173.9284000000007
$ws.Range("N138").Value = -21659.546

$ws.Range("H141").Value = 5205
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 5205
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 15615
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -25975

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8911.056
$ws.Range("I32").Value = 9258.764999999999
$ws.Range("K32").Value = 9258.764999999999
$ws.Range("M32").Value = -8971.764999999999

$ws.Range("H37").Value = 60937.312
$ws.Range("I37").Value = 14500
$ws.Range("J37").Value = 67571.21000000001
$ws.Range("K37").Value = 14500
$ws.Range("L37").Value = 67571.21000000001
$ws.Range("M37").Value = -14227
$ws.Range("N37").Value = -68117.21000000001

$ws.Range("H61").Value = 5690.8335
$ws.Range("I61").Value = 6136.615
$ws.Range("K61").Value = 6136.615
$ws.Range("M61").Value = -5924.615

$ws.Range("H122").Value = 3345.5588
$ws.Range("I122").Value = 2773
$ws.Range("J122").Value = 4542.727
$ws.Range("K122").Value = 8319
$ws.Range("L122").Value = 13628.181
$ws.Range("M122").Value = -5869
$ws.Range("N122").Value = -18528.181

$ws.Range("H132").Value = 3089.8447
$ws.Range("I132").Value = 2223.1765
$ws.Range("K132").Value = 6669.529500000001
$ws.Range("M132").Value = -4139.529500000001

$ws.Range("H136").Value = 5690.8335
$ws.Range("I136").Value = 6136.615
$ws.Range("K136").Value = 18409.845
$ws.Range("M136").Value = -15859.845

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 22974.068
$ws.Range("J99").Value = 12251.8
$ws.Range("L99").Value = 12251.8
$ws.Range("N99").Value = -15247.8

$ws.Range("H105").Value = 1521.2
$ws.Range("I105").Value = 1432
$ws.Range("J105").Value = 1967.2
$ws.Range("K105").Value = 1432
$ws.Range("L105").Value = 1967.2
$ws.Range("M105").Value = 315
$ws.Range("N105").Value = -5461.2

$ws.Range("H134").Value = 6505.8096
$ws.Range("I134").Value = 3602.077
$ws.Range("K134").Value = 10806.231
$ws.Range("M134").Value = -8271.231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4672.604
$ws.Range("I31").Value = 3225.7273
$ws.Range("K31").Value = 3225.7273
$ws.Range("M31").Value = -2930.7273

$ws.Range("H34").Value = 4672.604
$ws.Range("I34").Value = 3225.7273
$ws.Range("K34").Value = 3225.7273
$ws.Range("M34").Value = -3023.7273

$ws.Range("H94").Value = 1015.08
$ws.Range("I94").Value = 615.36365
$ws.Range("K94").Value = 615.36365
$ws.Range("M94").Value = -164.36365

$ws.Range("H99").Value = 10384
$ws.Range("I99").Value = 4937
$ws.Range("J99").Value = 15286.3
$ws.Range("K99").Value = 4937
$ws.Range("L99").Value = 15286.3
$ws.Range("M99").Value = -3439
$ws.Range("N99").Value = -18282.3

$ws.Range("H126").Value = 10384
$ws.Range("I126").Value = 4937
$ws.Range("J126").Value = 15286.3
$ws.Range("K126").Value = 14811
$ws.Range("L126").Value = 45858.89999999999
$ws.Range("M126").Value = -12341
$ws.Range("N126").Value = -50798.89999999999

$ws.Range("H134").Value = 1966.8695
$ws.Range("I134").Value = 1933.6316
$ws.Range("K134").Value = 5800.8948
$ws.Range("M134").Value = -3265.8948

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 31.444445
$ws.Range("J12").Value = 48.25
$ws.Range("L12").Value = 144.75
$ws.Range("N12").Value = -490.75

$ws.Range("H97").Value = 466.66666
$ws.Range("J97").Value = 466.66666
$ws.Range("L97").Value = 1399.99998
$ws.Range("N97").Value = -2391.99998

$ws.Range("H107").Value = 670.06665
$ws.Range("I107").Value = 681
$ws.Range("J107").Value = 666.0909
$ws.Range("K107").Value = 2043
$ws.Range("L107").Value = 1998.2727
$ws.Range("M107").Value = -123
$ws.Range("N107").Value = -5838.2727

$ws.Range("H113").Value = 1257.5714
$ws.Range("I113").Value = 1320
$ws.Range("J113").Value = 1222.8889
$ws.Range("K113").Value = 3960
$ws.Range("L113").Value = 3668.6667
$ws.Range("M113").Value = -1790
$ws.Range("N113").Value = -8008.6667

$ws.Range("H132").Value = 2665.5293
$ws.Range("J132").Value = 2663.7778
$ws.Range("L132").Value = 23974.0002
$ws.Range("N132").Value = -29034.0002

$ws.Range("H140").Value = 3496.3235
$ws.Range("I140").Value = 2220.0908
$ws.Range("K140").Value = 6660.2724
$ws.Range("M140").Value = -1480.2724

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5318.032
$ws.Range("I102").Value = 5594.9585
$ws.Range("K102").Value = 5594.9585
$ws.Range("M102").Value = -3972.9585

$ws.Range("H122").Value = 4148.5713
$ws.Range("J122").Value = 11494.5
$ws.Range("L122").Value = 34483.5
$ws.Range("N122").Value = -39383.5

$ws.Range("H126").Value = 3029.5454
$ws.Range("I126").Value = 2333.3333
$ws.Range("K126").Value = 6999.999899999999
$ws.Range("M126").Value = -4529.999899999999

$ws.Range("H132").Value = 23366.445
$ws.Range("I132").Value = 17649.834
$ws.Range("K132").Value = 52949.50199999999
$ws.Range("M132").Value = -50419.50199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3854.4
$ws.Range("J7").Value = 4800
$ws.Range("L7").Value = 4800
$ws.Range("N7").Value = -5024

$ws.Range("H100").Value = 7660.263
$ws.Range("I100").Value = 2267.2144
$ws.Range("J100").Value = 22760.8
$ws.Range("K100").Value = 2267.2144
$ws.Range("L100").Value = 22760.8
$ws.Range("M100").Value = -1726.2144
$ws.Range("N100").Value = -23842.8

$ws.Range("H126").Value = 3854.4
$ws.Range("J126").Value = 4800
$ws.Range("L126").Value = 14400
$ws.Range("N126").Value = -19340

$ws.Range("H132").Value = 3432.7827
$ws.Range("I132").Value = 3197.95
$ws.Range("J132").Value = 4998.3335
$ws.Range("K132").Value = 9593.849999999999
$ws.Range("L132").Value = 14995.0005
$ws.Range("M132").Value = -7063.849999999999
$ws.Range("N132").Value = -20055.0005

$ws.Range("H136").Value = 4133.9565
$ws.Range("I136").Value = 3713.4048
$ws.Range("K136").Value = 11140.2144
$ws.Range("M136").Value = -8590.214399999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2101.92
$ws.Range("I122").Value = 1661.2727
$ws.Range("K122").Value = 4983.8181
$ws.Range("M122").Value = -2533.8181

$ws.Range("H126").Value = 5011.5386
$ws.Range("I126").Value = 5011.5386
$ws.Range("K126").Value = 15034.6158
$ws.Range("M126").Value = -12564.6158

$ws.Range("H132").Value = 8773.950000000001
$ws.Range("I132").Value = 7926.7144
$ws.Range("J132").Value = 10750.833
$ws.Range("K132").Value = 23780.1432
$ws.Range("L132").Value = 32252.499
$ws.Range("M132").Value = -21250.1432
$ws.Range("N132").Value = -37312.499
